$wb = $excel.ActiveWorkbook

# Rename the "ThoracaicVertebrae" sheet to "ThoraicVertebrae"
$sheet = $wb.Worksheets.Item("ThoracaicVertebrae")
$sheet.Name = "ThoraicVertebrae"

# Update the label cell on the summary sheet that spells out the same name
$first = $wb.Worksheets.Item("Tissue-AM-Masses")
$first.Range("A9").Value = "ThoraicVertebrae"

# Activate the first sheet (Tissue-AM-Masses) and select A10
$first.Activate()
$first.Range("A10").Select()
